$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows above row 2 (shifts existing data/formulas down by 7 rows)
$ws.Rows("2:8").Insert()

# Fill in the 7 new rows of data (dates 44305 down to 44299)
$ws.Range("A7").Formula = "=(A8+1)"
$ws.Range("A6").Formula = "=(A7+1)"
$ws.Range("A5").Formula = "=(A6+1)"
$ws.Range("A4").Formula = "=(A5+1)"
$ws.Range("A3").Formula = "=(A4+1)"
$ws.Range("A2").Formula = "=(A3+1)"
$ws.Range("A8").Formula = "=(A9+1)"

$ws.Range("B2").Value = 0.75
$ws.Range("C2").Value = 1.21
$ws.Range("D2").Value = 0.98

$ws.Range("B3").Value = 0.77
$ws.Range("C3").Value = 1.19
$ws.Range("D3").Value = 0.98

$ws.Range("B4").Value = 0.79
$ws.Range("C4").Value = 1.18
$ws.Range("D4").Value = 0.99

$ws.Range("B5").Value = 0.81
$ws.Range("C5").Value = 1.17
$ws.Range("D5").Value = 0.99

$ws.Range("B6").Value = 0.83
$ws.Range("C6").Value = 1.1599999999999999
$ws.Range("D6").Value = 1

$ws.Range("B7").Value = 0.85
$ws.Range("C7").Value = 1.1599999999999999
$ws.Range("D7").Value = 1

$ws.Range("B8").Value = 0.87
$ws.Range("C8").Value = 1.1499999999999999
$ws.Range("D8").Value = 1.01

# Revise computed values for the following several rows (rolling averages recalculated
# now that new data has been added ahead of them)
$ws.Range("B9").Value = 0.91
$ws.Range("C9").Value = 1.1399999999999999
$ws.Range("D9").Value = 1.02

$ws.Range("B10").Value = 0.92
$ws.Range("C10").Value = 1.1299999999999999
$ws.Range("D10").Value = 1.03

$ws.Range("B11").Value = 0.94
$ws.Range("C11").Value = 1.1299999999999999
$ws.Range("D11").Value = 1.04

$ws.Range("B12").Value = 0.96
$ws.Range("C12").Value = 1.1299999999999999
$ws.Range("D12").Value = 1.05

$ws.Range("B13").Value = 0.98
$ws.Range("C13").Value = 1.1299999999999999
$ws.Range("D13").Value = 1.06

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 1.1399999999999999
$ws.Range("D14").Value = 1.06

$ws.Range("B15").Value = 1.01
$ws.Range("C15").Value = 1.1399999999999999
$ws.Range("D15").Value = 1.07

$ws.Range("B16").Value = 1.03
$ws.Range("C16").Value = 1.1499999999999999
$ws.Range("D16").Value = 1.08

$ws.Range("B17").Value = 1.04
$ws.Range("C17").Value = 1.1499999999999999
$ws.Range("D17").Value = 1.1000000000000001

$ws.Range("B18").Value = 1.05
$ws.Range("C18").Value = 1.1599999999999999
$ws.Range("D18").Value = 1.1000000000000001

$ws.Range("B19").Value = 1.06
$ws.Range("C19").Value = 1.17
$ws.Range("D19").Value = 1.1100000000000001

$ws.Range("B20").Value = 1.08
$ws.Range("C20").Value = 1.18
$ws.Range("D20").Value = 1.1299999999999999

$ws.Range("B21").Value = 1.0900000000000001
$ws.Range("C21").Value = 1.18
$ws.Range("D21").Value = 1.1299999999999999

$ws.Range("B22").Value = 1.1000000000000001
$ws.Range("C22").Value = 1.19
$ws.Range("D22").Value = 1.1399999999999999

$ws.Range("B23").Value = 1.1000000000000001
$ws.Range("C23").Value = 1.2
$ws.Range("D23").Value = 1.1499999999999999

$ws.Range("B24").Value = 1.1100000000000001
$ws.Range("C24").Value = 1.21
$ws.Range("D24").Value = 1.1599999999999999

$ws.Range("B25").Value = 1.1200000000000001
$ws.Range("C25").Value = 1.23
$ws.Range("D25").Value = 1.17

$ws.Range("B26").Value = 1.1299999999999999
$ws.Range("C26").Value = 1.23
$ws.Range("D26").Value = 1.18

$ws.Range("B27").Value = 1.1299999999999999
$ws.Range("C27").Value = 1.24
$ws.Range("D27").Value = 1.19

$ws.Range("B28").Value = 1.1399999999999999
$ws.Range("C28").Value = 1.25
$ws.Range("D28").Value = 1.2

# Update view: selection
$ws.Range("M26:M27").Select()
